$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Fill in the data rows for the existing Table8 (D3:E5)
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = "A"
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = "B"

# Add data + new table (Table7) at G3:G4
$ws.Range("G3").Value = "Column1"
$ws.Range("G4").Value = "aaa"
$lo = $ws.ListObjects.Add(1, $ws.Range("G3:G4"), 0, 1)

# Column width tweaks matching surrounding columns
$ws.Columns.Item(7).ColumnWidth = 9.59
$ws.Columns.Item(9).ColumnWidth = 9.59

# Final selection
[void]$ws.Range("I8").Select()
